$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 86.28570999999999
$ws.Range("I11").Value = 86.28570999999999
$ws.Range("K11").Value = 86.28570999999999
$ws.Range("M11").Value = 53.71429000000001

$ws.Range("H17").Value = 3166.6667
$ws.Range("I17").Value = 1500
$ws.Range("J17").Value = 4000
$ws.Range("K17").Value = 4500
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = -4332
$ws.Range("N17").Value = -12336

$ws.Range("H40").Value = 3136.8
$ws.Range("I40").Value = 1874.75
$ws.Range("K40").Value = 1874.75
$ws.Range("M40").Value = -1699.75

$ws.Range("H51").Value = 4997
$ws.Range("J51").Value = 4997
$ws.Range("L51").Value = 4997
$ws.Range("N51").Value = -5965

$ws.Range("H62").Value = 4833.4707
$ws.Range("J62").Value = 3450
$ws.Range("L62").Value = 3450
$ws.Range("N62").Value = -4698

$ws.Range("H65").Value = 4833.4707
$ws.Range("J65").Value = 3450
$ws.Range("L65").Value = 17250
$ws.Range("N65").Value = -23490

$ws.Range("H116").Value = 6319.4443
$ws.Range("I116").Value = 5990.3335
$ws.Range("K116").Value = 5990.3335
$ws.Range("M116").Value = -2548.3335

$ws.Range("H132").Value = 2941.5366
$ws.Range("I132").Value = 2886.6572
$ws.Range("K132").Value = 8659.971600000001
$ws.Range("M132").Value = -6129.971600000001

$ws.Range("H138").Value = 2692
$ws.Range("I138").Value = 1096
$ws.Range("J138").Value = 3490
$ws.Range("K138").Value = 3288
$ws.Range("L138").Value = 10470
$ws.Range("M138").Value = 1852
$ws.Range("N138").Value = -20750

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6274.9556
$ws.Range("I32").Value = 5281.2046
$ws.Range("K32").Value = 5281.2046
$ws.Range("M32").Value = -4994.2046

$ws.Range("H37").Value = 22500
$ws.Range("I37").Value = 17500
$ws.Range("K37").Value = 17500
$ws.Range("M37").Value = -17227

$ws.Range("H46").Value = 4743.3335
$ws.Range("I46").Value = 4990
$ws.Range("J46").Value = 4496.6665
$ws.Range("K46").Value = 4990
$ws.Range("L46").Value = 4496.6665
$ws.Range("M46").Value = -4671
$ws.Range("N46").Value = -5134.6665

$ws.Range("H80").Value = 38000

$ws.Range("H83").Value = 38000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3671
$ws.Range("I31").Value = 2691.889
$ws.Range("K31").Value = 2691.889
$ws.Range("M31").Value = -2396.889

$ws.Range("H34").Value = 3671
$ws.Range("I34").Value = 2691.889
$ws.Range("K34").Value = 2691.889
$ws.Range("M34").Value = -2489.889

$ws.Range("H58").Value = 5228.375
$ws.Range("I58").Value = 6438.8335
$ws.Range("K58").Value = 6438.8335
$ws.Range("M58").Value = -6235.8335

$ws.Range("H60").Value = 25000
$ws.Range("J60").Value = 25000
$ws.Range("L60").Value = 25000
$ws.Range("N60").Value = -26022

$ws.Range("H64").Value = 45000
$ws.Range("J64").Value = 45000
$ws.Range("L64").Value = 45000
$ws.Range("N64").Value = -45496

$ws.Range("H67").Value = 45000
$ws.Range("J67").Value = 45000
$ws.Range("L67").Value = 45000
$ws.Range("N67").Value = -46716

$ws.Range("H134").Value = 1393.3846
$ws.Range("I134").Value = 1390.625
$ws.Range("K134").Value = 4171.875
$ws.Range("M134").Value = -1636.875

$ws.Range("H136").Value = 5228.375
$ws.Range("I136").Value = 6438.8335
$ws.Range("K136").Value = 19316.5005
$ws.Range("M136").Value = -16766.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 545.2
$ws.Range("I5").Value = 516.125
$ws.Range("J5").Value = 661.5
$ws.Range("K5").Value = 1548.375
$ws.Range("L5").Value = 1984.5
$ws.Range("M5").Value = -1436.375
$ws.Range("N5").Value = -2208.5

$ws.Range("H14").Value = 42379.777
$ws.Range("I14").Value = 42379.777
$ws.Range("K14").Value = 127139.331
$ws.Range("M14").Value = -126966.331

$ws.Range("H23").Value = 588.4
$ws.Range("J23").Value = 487.1111
$ws.Range("L23").Value = 1461.3333
$ws.Range("N23").Value = -1931.3333

$ws.Range("H33").Value = 109.85714
$ws.Range("I33").Value = 94.833336
$ws.Range("K33").Value = 569.000016
$ws.Range("M33").Value = -286.000016

$ws.Range("H81").Value = 8839
$ws.Range("J81").Value = 8839
$ws.Range("L81").Value = 26517
$ws.Range("N81").Value = -28763

$ws.Range("H84").Value = 8839
$ws.Range("J84").Value = 8839
$ws.Range("L84").Value = 79551
$ws.Range("N84").Value = -90783

$ws.Range("H135").Value = 545.2
$ws.Range("I135").Value = 516.125
$ws.Range("J135").Value = 661.5
$ws.Range("K135").Value = 4645.125
$ws.Range("L135").Value = 5953.5
$ws.Range("M135").Value = -2110.125
$ws.Range("N135").Value = -11023.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4072.2727
$ws.Range("I80").Value = 2198.75
$ws.Range("J80").Value = 5142.857
$ws.Range("K80").Value = 2198.75
$ws.Range("L80").Value = 5142.857
$ws.Range("M80").Value = -1200.75
$ws.Range("N80").Value = -7138.857

$ws.Range("H83").Value = 4072.2727
$ws.Range("I83").Value = 2198.75
$ws.Range("J83").Value = 5142.857
$ws.Range("K83").Value = 10993.75
$ws.Range("L83").Value = 25714.285
$ws.Range("M83").Value = -6001.75
$ws.Range("N83").Value = -35698.285

$ws.Range("H97").Value = 704.1818
$ws.Range("I97").Value = 752.2632
$ws.Range("J97").Value = 399.66666
$ws.Range("K97").Value = 752.2632
$ws.Range("L97").Value = 399.66666
$ws.Range("M97").Value = -256.2632
$ws.Range("N97").Value = -1391.66666

$ws.Range("H132").Value = 995
$ws.Range("I132").Value = 995
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2985
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -455
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1425
$ws.Range("I55").Value = 1916.6666
$ws.Range("J55").Value = 933.3333
$ws.Range("K55").Value = 1916.6666
$ws.Range("L55").Value = 933.3333
$ws.Range("M55").Value = -1743.6666
$ws.Range("N55").Value = -1279.3333

$ws.Range("H100").Value = 1957.6364
$ws.Range("I100").Value = 2003.8889
$ws.Range("J100").Value = 1749.5
$ws.Range("K100").Value = 2003.8889
$ws.Range("L100").Value = 1749.5
$ws.Range("M100").Value = -1462.8889
$ws.Range("N100").Value = -2831.5

$ws.Range("H132").Value = 17055.705
$ws.Range("I132").Value = 16871.688
$ws.Range("J132").Value = 20000
$ws.Range("K132").Value = 50615.064
$ws.Range("L132").Value = 60000
$ws.Range("M132").Value = -48085.064
$ws.Range("N132").Value = -65060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 32999.668
$ws.Range("I51").Value = 28999.5
$ws.Range("K51").Value = 28999.5
$ws.Range("M51").Value = -28489.5

$ws.Range("H100").Value = 822.1
$ws.Range("I100").Value = 789.1429000000001
$ws.Range("J100").Value = 899
$ws.Range("K100").Value = 1578.2858
$ws.Range("L100").Value = 1798
$ws.Range("M100").Value = -1037.2858
$ws.Range("N100").Value = -2880

$ws.Range("H113").Value = 1575.2941
$ws.Range("I113").Value = 677.0909
$ws.Range("J113").Value = 3222
$ws.Range("K113").Value = 2031.2727
$ws.Range("L113").Value = 9666
$ws.Range("M113").Value = 138.7273
$ws.Range("N113").Value = -14006

$ws.Range("H122").Value = 1281.9
$ws.Range("I122").Value = 1281.9
$ws.Range("K122").Value = 3845.7
$ws.Range("M122").Value = -1395.7

$ws.Range("H126").Value = 1543.3334
$ws.Range("I126").Value = 1623.75
$ws.Range("K126").Value = 4871.25
$ws.Range("M126").Value = -2401.25

$ws.Range("H132").Value = 441.6
$ws.Range("I132").Value = 377.5
$ws.Range("K132").Value = 1132.5
$ws.Range("M132").Value = 1397.5

$ws.Range("H136").Value = 3028.35
$ws.Range("I136").Value = 3103.9412
$ws.Range("J136").Value = 2600
$ws.Range("K136").Value = 9311.8236
$ws.Range("L136").Value = 7800
$ws.Range("M136").Value = -6761.8236
$ws.Range("N136").Value = -12900
